$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.699.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.28%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.076.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.84%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.91%  '

# Row 6
$ws.Range("E6").Value = '  -0.56%  '

# Row 7
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.25'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.35%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.390'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.51%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0783'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.25%  '

# Row 11
$ws.Range("E11").Value = '  +2.54%  '

# Row 12
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.385.88'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.88%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.70%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.79%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.770'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.34%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.044.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.40%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.586.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.56%  '

# Row 19
$ws.Range("E19").Value = '  -0.34%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.09%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0831'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("E23").Value = '  -0.11%  '

# Row 24
$ws.Range("E24").Value = '  -0.54%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.99%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.29%  '

# Row 27
$ws.Range("E27").Value = '  +2.53%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.39%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.32%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.23%  '

# Row 31
$ws.Range("E31").Value = '  +1.87%  '

# Row 32
$ws.Range("E32").Value = '  +0.51%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0628'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.01%  '

# Row 34
$ws.Range("E34").Value = '  +1.14%  '

# Row 35
$ws.Range("E35").Value = '  -3.81%  '

# Row 36
$ws.Range("E36").Value = '  +2.79%  '

# Row 37
$ws.Range("E37").Value = '  -3.28%  '

# Row 38
$ws.Range("E38").Value = '  -0.13%  '

# Row 39
$ws.Range("E39").Value = '  -5.98%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0976'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.05%  '

# Row 42
$ws.Range("E42").Value = '  +0.62%  '

# Row 43
$ws.Range("E43").Value = '  -2.81%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.450.05'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.82%  '

# Row 45
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.99%  '

# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.25%  '

# Row 47
$ws.Range("E47").Value = '  +0.84%  '

# Row 48
$ws.Range("E48").Value = '  +0.42%  '

# Row 49
$ws.Range("E49").Value = '  +0.25%  '

# Row 50
$ws.Range("E50").Value = '  -1.47%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.267.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.08%  '
